# Append a new row to the project-timeline table documenting the
# "biomes system" work completed on 19/12/2021.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Rows.Add() with no argument appends a new row at the end of the table,
# inheriting the cell formatting (shading/borders/widths) of the last row -
# exactly what we want here, since the new row mirrors the row above it.
$newRow = $table.Rows.Add()
$r = $newRow.Index

$table.Cell($r, 1).Range.Text = "19/12/2021"
$table.Cell($r, 2).Range.Text = "2 Hours 10 Minutes"
$table.Cell($r, 3).Range.Text = "World Generation " + [char]0x2013 + " Objective 1"
$table.Cell($r, 4).Range.Text = "Implemented the biome system properly " + [char]0x2013 + " allowing for the generation of a world with different climates. As of this stage, there should be sufficient progress to move onto the next goals as well as update the dissertation document."
